# Generate Report for Handback
#
# Refresh the localization status report with the results of the latest
# handback run: the zh-cn/de-de files came back in sync with en-US, so the
# status message changes, the "Latest Handback DateTime" timestamps move
# forward, and the stale "version mismatch" error notes are cleared.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: per-language status cells.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# zh-cn detail sheet.
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-20 08:56:44"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P2").Style = "Normal"

# de-de detail sheet.
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-20 08:56:50"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P2").Style = "Normal"

# Widen the Status columns to comfortably fit the new, longer status text,
# and shrink the now-empty Error Detail columns.
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 29.9777047293527
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 29.9777047293527
$wsZhCn.Range("P1").EntireColumn.ColumnWidth = 13.7470528738839
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 29.9777047293527
$wsDeDe.Range("P1").EntireColumn.ColumnWidth = 13.7470528738839
